$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Edge list: column A = "from", column B = "to", column C = "length".
# Add the charging-node edges for tugs 107/108/109/110 right after the
# existing data (previously ending at row 253).
$newEdges = @(
    @(108, 109, 0.5),
    @(109, 108, 0.5),
    @(107, 110, 0.5),
    @(110, 107, 0.5)
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newEdges.Count; $i++) {
    $row = $startRow + $i
    $edge = $newEdges[$i]
    $ws.Cells.Item($row, 1).Value = $edge[0]
    $ws.Cells.Item($row, 2).Value = $edge[1]
    $ws.Cells.Item($row, 3).Value = $edge[2]
}

# Match the author's final view state: scrolled near the bottom of the
# sheet with A102 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 241
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("A102").Select()
